$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Chin Pok%CHAN%chinpok0301@gmail.com%1,         Ngai Sze%WONG%candy_wong@cuhk.edu.hk%2,         Ngai Sze%WONG%candy_wong@cuhk.edu.hk%0,         Chi Chiu%LEUNG%ccleungpnc@netvigator.com%2,         Chi Chiu%LEUNG%ccleungpnc@netvigator.com%0,         Shui Shan%LEE%sslee@cuhk.edu.hk%1]'
$ws.Range("E3").Value = '[Chan%Crystal Ying%coreGivesNoEmail%1,        Chan%Ka%coreGivesNoEmail%1,        Ho%Pak-Leung%coreGivesNoEmail%1,        Lam%Kin Bong Hubert%coreGivesNoEmail%1,        Lee%Pak-Wing%coreGivesNoEmail%1]'
$ws.Range("E4").Value = '[Khee-Siang%Chan%NULL%1,         Fu-Wen%Liang%NULL%1,         Hung-Jen%Tang%NULL%1,         Han Siong%Toh%NULL%1,         Wen-Liang%Yu%NULL%1]'
$ws.Range("E5").Value = '[Young June%Choe%NULL%1,         Jong-Koo%Lee%NULL%2,         Jong-Koo%Lee%NULL%0]'
$ws.Range("E6").Value = '[Ermengol%Coma Redon%NULL%1,         Nuria%Mora%NULL%1,         Albert%Prats-Uribe%NULL%1,         Francesc%Fina Avilés%NULL%2,         Francesc%Fina Avilés%NULL%0,         Daniel%Prieto-Alhambra%NULL%1,         Manuel%Medina%NULL%2,         Manuel%Medina%NULL%0]'
$ws.Range("E7").Value = '[Benjamin J%Cowling%NULL%2,         Sheikh Taslim%Ali%NULL%2,         Tiffany W Y%Ng%NULL%2,         Tim K%Tsang%NULL%2,         Julian C M%Li%NULL%2,         Min Whui%Fong%NULL%2,         Qiuyan%Liao%NULL%2,         Mike YW%Kwan%NULL%2,         So Lun%Lee%NULL%2,         Susan S%Chiu%NULL%2,         Joseph T%Wu%NULL%2,         Peng%Wu%pengwu@hku.hk%2,         Gabriel M%Leung%NULL%2]'
$ws.Range("E8").Value = '[Takahiro%Itaya%NULL%1,         Yuki%Furuse%NULL%1,         Kazuaki%Jindai%NULL%1]'
$ws.Range("E9").Value = '[Xiangsha%Kong%NULL%1,         Feng%Liu%NULL%0,         Haibo%Wang%NULL%1,         Ruifeng%Yang%NULL%1,         Dongbo%Chen%NULL%1,         Xiaoxiao%Wang%NULL%1,         Fengmin%Lu%lu.fengmin@hsc.pku.edu.cn%1,         Huiying%Rao%raohuiying@pkuph.edu.cn%1,         Hongsong%Chen%chenhongsong@bjmu.edu.cn%1]'
$ws.Range("E10").Value = '[Shu-Chen%Kuo%NULL%1,         Shu-Man%Shih%NULL%1,         Li-Hsin%Chien%NULL%1,         Chao A.%Hsiung%NULL%1]'
$ws.Range("E11").Value = '[Hyunju%Lee%NULL%1,         Heeyoung%Lee%NULL%1,         Kyoung-Ho%Song%NULL%0,         Eu Suk%Kim%NULL%2,         Eu Suk%Kim%NULL%0,         Jeong Su%Park%NULL%1,         Jongtak%Jung%NULL%1,         Soyeon%Ahn%NULL%1,         Eun Kyeong%Jeong%NULL%1,         Hyekyung%Park%NULL%1,         Hong Bin%Kim%hbkimmd@snu.ac.kr%0]'
$ws.Range("E12").Value = '[Ji Yun%Noh%NULL%0,         Hye%Seong%NULL%0,         Hye%Seong%NULL%0,         Jin Gu%Yoon%NULL%0,         Jin Gu%Yoon%NULL%0,         Joon Young%Song%NULL%0,         Joon Young%Song%NULL%0,         Hee Jin%Cheong%NULL%0,         Hee Jin%Cheong%NULL%0,         Woo Joo%Kim%NULL%0,         Woo Joo%Kim%NULL%0]'
$ws.Range("E13").Value = '[Di%Wu%NULL%1,         Jianyun%Lu%NULL%1,         Yanhui%Liu%NULL%1,         Zhoubin%Zhang%NULL%1,         Lei%Luo%NULL%1]'
$ws.Range("E14").Value = '[Chih-Jen%Yang%NULL%1,         Tun-Chieh%Chen%NULL%2,         Tun-Chieh%Chen%NULL%0,         Shin-Huei%Kuo%NULL%1,         Min-Han%Hsieh%NULL%1,         Yen-Hsu%Chen%NULL%1]'
$ws.Range("E15").Value = '[ G. M.%Young%null%1,          X.% Peng%null%1,          A.% Rebaza%null%1,          S.% Bermejo%null%1,          D.% Chang%null%1,          L.% Sharma%null%1,          C. % Dela Cruz%null%1]'
